$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Log new hours entries (Mr. Pink sheet values)
$ws.Range("E6").Value = 41062
$ws.Range("F6").Value = 6.5

$ws.Range("E7").Value = 41063
$ws.Range("F7").Value = 2

$ws.Range("B8").Value = 41062
$ws.Range("C8").Value = 6.5

# Update the selected cell to reflect where the user left off editing
$ws.Range("R12").Select()
